$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the "Upper Left Cell" / "Lower Right Cell" index references in column D
# (rows 5-11) for the Scenario/Directory/DSS Path/Date index blocks from row 25
# to row 26, reflecting the updated init file layout for scenario 46.
$ws.Range("D5").Value = "A26"
$ws.Range("D6").Value = "B26"
$ws.Range("D7").Value = "C26"
$ws.Range("D8").Value = "G26"
$ws.Range("D9").Value = "H26"
$ws.Range("D10").Value = "I26"
$ws.Range("D11").Value = "J26"
